$wb = $excel.ActiveWorkbook

# 1. Add a new worksheet "Plan1" at the end of the workbook.
$sheet1 = $wb.Worksheets.Item("Baby Step 01")
$sheet2 = $wb.Worksheets.Item("Baby Step 02")
$plan1 = $wb.Worksheets.Add($null, $sheet2)
$plan1.Name = "Plan1"

# 2. Populate "Plan1" with the example/leftover tables from the modelling session.
$plan1.Range("E1").Value = "0-1"

$plan1.Range("A2").Value = "label 1"
$plan1.Range("B2").Value = 0.5
$plan1.Range("C2").Value = 65
$plan1.Range("D2").Value = 90

$plan1.Range("A3").Value = "label 2"
$plan1.Range("B3").Value = 4.5
$plan1.Range("C3").Value = 34
$plan1.Range("D3").Value = 0.5

$plan1.Range("I8").Value = "Coisa 1"
$plan1.Range("J8").Value = "adj 1"
$plan1.Range("K8").Value = "adj 2"
$plan1.Range("L8").Value = "adj 3"
$plan1.Range("M8").Value = "resultado"

$plan1.Range("H9").Value = " "
$plan1.Range("I9").Value = "Coisa 2"
$plan1.Range("J9").Value = "adj 1"
$plan1.Range("K9").Value = "adj 2"
$plan1.Range("L9").Value = "adj 3"
$plan1.Range("M9").Value = "resultado"

$plan1.Range("I10").Value = "Coisa 3"
$plan1.Range("J10").Value = "adj 1"
$plan1.Range("K10").Value = "adj 2"
$plan1.Range("L10").Value = "adj 3"
$plan1.Range("M10").Value = "extrapolar"

$plan1.Range("I12").Value = "Renda"
$plan1.Range("J12").Value = "Escolaridade"
$plan1.Range("K12").Value = "Região"
$plan1.Range("L12").Value = "Probabilidade de ser preso 0-1"

$plan1.Range("H13").Value = "pessoa 1"
$plan1.Range("I13").Value = "baixa"
$plan1.Range("J13").Value = "médio"
$plan1.Range("K13").Value = "sudeste"
$plan1.Range("L13").Value = 0.2

$plan1.Range("H14").Value = "pessoa 2"
$plan1.Range("I14").Value = "alta"
$plan1.Range("J14").Value = "superior"
$plan1.Range("K14").Value = "sul"
$plan1.Range("L14").Value = 0.1

$plan1.Range("H15").Value = "pessoa 3"
$plan1.Range("I15").Value = "baixa"
$plan1.Range("J15").Value = "fundamental"
$plan1.Range("K15").Value = "norte"
$plan1.Range("L15").Value = 0.4
